$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Write the new cell text/number values in the exact order the strings
#    were originally appended to the shared-strings table (column-major
#    within each translation block: English column, then filename, then
#    Russian column, then transliterated column) so that new shared strings
#    land at the same indices as the target workbook.
# ---------------------------------------------------------------------------

$ws.Range("C49").Value = ' The horrible enemies within\n[CS:P]Brine Cave[CR]...[K] What do you suppose they are?'
$ws.Range("C50").Value = ' ...[K]Eek! Thinking about them\nfrightens me!'
$ws.Range("A49").Value = 'SCRIPT/T01P02A/um2201.ssb'
$ws.Range("D49").Value = ' Ужасные враги в [CS:P]Пещере у Моря[CR]...[K]\nКак думаете, кто они?'
$ws.Range("D50").Value = ' ...[K]Иии! Даже одна мысль о них\nпугает меня!'
$ws.Range("E49").Value = ' Ôçàòîúå âñàãé â [CS:P]Ðåþåñå ô Íïñÿ[CR]...[K]\nËàë äôíàåóå, ëóï ïîé?'
$ws.Range("E50").Value = ' ...[K]Ééé! Äàçå ïäîà íúòìû ï îéö\nðôãàåó íåîÿ!'
$ws.Range("C51").Value = ' It will be all right![K]\nIt\''s [partner] and [hero]!'
$ws.Range("C52").Value = ' I\''m certain that you will be able\nto stop the planet\''s paralysis!'
$ws.Range("A51").Value = 'SCRIPT/T01P02A/um2408.ssb'
$ws.Range("D51").Value = ' Всё будет хорошо![K]\nВы же [partner] и [hero]!'
$ws.Range("D52").Value = ' Я уверена, что вам по силам\nостановить планетарный паралич!'
$ws.Range("E51").Value = ' Âòæ áôäåó öïñïšï![K]\nÂú çå [partner] é [hero]!'
$ws.Range("E52").Value = ' Ÿ ôâåñåîà, œóï âàí ðï òéìàí\nïòóàîïâéóû ðìàîåóàñîúê ðàñàìéœ!'

# Numeric "Line number" column
$ws.Range("B49").Value = 395
$ws.Range("B50").Value = 398
$ws.Range("B51").Value = 373
$ws.Range("B52").Value = 376

# ---------------------------------------------------------------------------
# 2) Formatting. Reuse existing cell styles via copy/paste-formats so no new
#    style entries are introduced into the stylesheet.
# ---------------------------------------------------------------------------

# Row 48 becomes the closing (bottom-bordered) row of its block, matching the
# look of every other block end (e.g. row 46). Gains an empty A48 cell too.
$ws.Range("A46:E46").Copy()
$ws.Range("A48:E48").PasteSpecial(-4122)

# Row 49: first row of a new block -> "name" style (plain, with filename in A)
$ws.Range("A47:E47").Copy()
$ws.Range("A49:E49").PasteSpecial(-4122)

# Row 50: closing row of that block -> bottom-bordered style
$ws.Range("A46:E46").Copy()
$ws.Range("A50:E50").PasteSpecial(-4122)

# Row 51: first row of the next block -> "name" style again
$ws.Range("A47:E47").Copy()
$ws.Range("A51:E51").PasteSpecial(-4122)

# Row 52: continuation row (no filename, and no A-cell at all) -> plain style
$ws.Range("B44:E44").Copy()
$ws.Range("B52:E52").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Row heights for the four brand-new rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(49).RowHeight = 43.2
$ws.Rows.Item(50).RowHeight = 21.6
$ws.Rows.Item(51).RowHeight = 43.2
$ws.Rows.Item(52).RowHeight = 31.8

# ---------------------------------------------------------------------------
# 4) Selection / scroll position, matching what Excel leaves behind after the
#    edit (selection lands on the last edited cell).
# ---------------------------------------------------------------------------
$ws.Range("E52").Select() | Out-Null
